$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update error-table statistics for rows 2-8 (columns B:G) to reflect the
# corrected selection scopes used when computing these aggregates.
$ws.Range("B2").Value = -0.1243020608379911
$ws.Range("C2").Value = 2.131588938050424
$ws.Range("D2").Value = 21.22560083334914
$ws.Range("E2").Value = 4.607125007349935
$ws.Range("F2").Value = 4.713826149300701
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = -0.246485675231136
$ws.Range("C3").Value = 2.345373773945208
$ws.Range("D3").Value = 19.33065173826034
$ws.Range("E3").Value = 4.396663705386204
$ws.Range("F3").Value = 4.498154218974101
$ws.Range("G3").Value = 21

$ws.Range("B4").Value = -0.6924834931754604
$ws.Range("C4").Value = 1.814664499321882
$ws.Range("D4").Value = 10.78567221865277
$ws.Range("E4").Value = 3.284154719049146
$ws.Range("F4").Value = 3.29371631896608
$ws.Range("G4").Value = 20

$ws.Range("B5").Value = -0.2778078423057671
$ws.Range("C5").Value = 1.655704769799213
$ws.Range("D5").Value = 12.12863580349151
$ws.Range("E5").Value = 3.482619101120809
$ws.Range("F5").Value = 3.56664885913509
$ws.Range("G5").Value = 19

$ws.Range("B6").Value = -0.2358027714400601
$ws.Range("C6").Value = 1.759610300011659
$ws.Range("D6").Value = 10.62893045379797
$ws.Range("E6").Value = 3.260204050944967
$ws.Range("F6").Value = 3.345936034707538
$ws.Range("G6").Value = 18

$ws.Range("B7").Value = -0.298460077026763
$ws.Range("C7").Value = 1.862981289594231
$ws.Range("D7").Value = 11.19011793689163
$ws.Range("E7").Value = 3.345163364753899
$ws.Range("F7").Value = 3.434363767755224
$ws.Range("G7").Value = 17

$ws.Range("B8").Value = -0.1522653927570955
$ws.Range("C8").Value = 1.879133647755224
$ws.Range("D8").Value = 12.22827700966262
$ws.Range("E8").Value = 3.496895338677241
$ws.Range("F8").Value = 3.608152585070817
$ws.Range("G8").Value = 16
